$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.617.82"
$ws.Range("E2").Value = "'  -0.94%  "
$ws.Range("D3").Value = "'2.277.75"
$ws.Range("E3").Value = "'  -0.97%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'304.52"
$ws.Range("E5").Value = "'  +1.49%  "
$ws.Range("D6").Value = "'96.21"
$ws.Range("E6").Value = "'  -1.66%  "
$ws.Range("D7").Value = "'0.506"
$ws.Range("E9").Value = "'  -3.12%  "
$ws.Range("D10").Value = "'35.52"
$ws.Range("E10").Value = "'  -1.53%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "'  -0.18%  "
$ws.Range("D12").Value = "'18.18"
$ws.Range("E12").Value = "'  +2.52%  "
$ws.Range("E13").Value = "'  +0.73%  "
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("D15").Value = "'2.633.66"
$ws.Range("D16").Value = "'2.281.84"
$ws.Range("E16").Value = "'  -1.28%  "
$ws.Range("D17").Value = "'0.776"
$ws.Range("E17").Value = "'  -1.45%  "
$ws.Range("D18").Value = "'42.567.66"
$ws.Range("E18").Value = "'  -0.79%  "
$ws.Range("D19").Value = "'12.97"
$ws.Range("E19").Value = "'  +1.61%  "
$ws.Range("D20").Value = "'0.0₃0892"
$ws.Range("E20").Value = "'  -1.98%  "
$ws.Range("E21").Value = "'  -2.21%  "
$ws.Range("D22").Value = "'67.11"
$ws.Range("E22").Value = "'  -1.66%  "
$ws.Range("D23").Value = "'235.41"
$ws.Range("E24").Value = "'  -2.13%  "
$ws.Range("E25").Value = "'  +0.05%  "
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "'  +0.66%  "
$ws.Range("E27").Value = "'  +0.01%  "
$ws.Range("D28").Value = "'25.07"
$ws.Range("E28").Value = "'  +0.54%  "
$ws.Range("D29").Value = "'166.09"
$ws.Range("E29").Value = "'  +1.60%  "
$ws.Range("E30").Value = "'  +0.71%  "
$ws.Range("D31").Value = "'9.03"
$ws.Range("E31").Value = "'  -1.15%  "
$ws.Range("D32").Value = "'33.01"
$ws.Range("E32").Value = "'  +0.02%  "
$ws.Range("E33").Value = "'  +0.08%  "
$ws.Range("E34").Value = "'  -0.67%  "
$ws.Range("E35").Value = "'  -3.13%  "
$ws.Range("D36").Value = "'17.56"
$ws.Range("E36").Value = "'  -2.65%  "
$ws.Range("E37").Value = "'  -1.13%  "
$ws.Range("D38").Value = "'0.0689"
$ws.Range("E38").Value = "'  -0.99%  "
$ws.Range("E39").Value = "'  -0.97%  "
$ws.Range("E40").Value = "'  -2.23%  "
$ws.Range("E41").Value = "'  -1.57%  "
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "'  -3.59%  "
$ws.Range("D43").Value = "'2.001.12"
$ws.Range("E43").Value = "'  -0.96%  "
$ws.Range("E44").Value = "'  -2.86%  "
$ws.Range("D45").Value = "'18.08"
$ws.Range("E45").Value = "'  +3.94%  "
$ws.Range("D46").Value = "'9.96"
$ws.Range("E46").Value = "'  -3.97%  "
$ws.Range("D47").Value = "'2.07"
$ws.Range("E47").Value = "'  -8.37%  "
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "'  -2.51%  "
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "'  +5.16%  "
$ws.Range("D50").Value = "'53.50"
$ws.Range("E50").Value = "'  -1.45%  "
$ws.Range("D51").Value = "'2.501.67"
$ws.Range("E51").Value = "'  -0.95%  "
